$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts PRESIDENT.. right by one)
$ws.Columns("C").Insert()

# Give the new header cell the same look as the other header cells, then set its text
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value2 = "ACADEMIC_YEAR"

# Remove the sample data row (old row 2)
$ws.Rows(2).Delete()

# Column C keeps a plain custom width of 15 (not best-fit, since it's a brand new column)
$ws.Columns("C").ColumnWidth = 14.17

# Update active selection to B6
$ws.Range("B6").Select() | Out-Null
